# Update cryptocurrency price (D) and volume/1h (E) values per latest symbol list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.40%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'40.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.97%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.129"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.09%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07615"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.23%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.622"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.33%"
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'0.42%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.9015"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'2.53%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1097"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'9.90%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1758"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.65%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.09116"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.19%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.04183"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-5.09%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.1051"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.36%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.001260"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.13%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.005832"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.82%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'3.354"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.06%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'4.258"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.42%"
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'-0.91%"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'6.592"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-6.33%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1366"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'2.02%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D22").Value = "'0.04065"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.29%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.001223"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.94%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.004087"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.89%"
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'6.58%"
$ws.Range("E25").Style = "Normal"

$ws.Range("E38").Value = "'1.50%"
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.05174"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.46%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.007782"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.64%"
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'-1.79%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.006765"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'6.97%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.001950"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-1.49%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.007967"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-6.15%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.3333"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'8.86%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00007022"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'8.09%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.18%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.03239"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'375.94%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.004207"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-39.95%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.00002104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.18%"
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0002003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.18%"
$ws.Range("E51").Style = "Normal"
